$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 18:25"

# Update country names whose sort position shifted (column A)
$ws.Range("A37").Value = "Irak"
$ws.Range("A38").Value = "Ucrania"
$ws.Range("A39").Value = "Portugal"
$ws.Range("A57").Value = "Ghana"
$ws.Range("A58").Value = "Moldavia"
$ws.Range("A59").Value = "Guatemala"
$ws.Range("A63").Value = "Argelia"
$ws.Range("A64").Value = "Dinamarca"
$ws.Range("A65").Value = "Corea del Sur"
$ws.Range("A66").Value = "Camerun"
$ws.Range("A127").Value = "Jordania"
$ws.Range("A128").Value = "Congo"
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A201").Value = "Laos"
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A212").Value = "Seychelles"
$ws.Range("A213").Value = "Montserrat"

# Update refreshed statistics for affected countries (columns B-H)
$ws.Range("B4").Value = 2512794
$ws.Range("C4").Value = 8206
$ws.Range("D4").Value = 1052880
$ws.Range("E4").Value = 1332995
$ws.Range("G4").Value = 139
$ws.Range("H4").Value = 126919

$ws.Range("B7").Value = 506972
$ws.Range("C7").Value = 15802
$ws.Range("D7").Value = 294988
$ws.Range("E7").Value = 196322
$ws.Range("G7").Value = 354
$ws.Range("H7").Value = 15662

$ws.Range("B11").Value = 263360
$ws.Range("C11").Value = 4296
$ws.Range("D11").Value = 223431
$ws.Range("E11").Value = 34861
$ws.Range("G11").Value = 165
$ws.Range("H11").Value = 5068

$ws.Range("B12").Value = 239961
$ws.Range("C12").Value = 255
$ws.Range("D12").Value = 187615
$ws.Range("E12").Value = 17638
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 34708

$ws.Range("B16").Value = 194013
$ws.Range("C16").Value = 228
$ws.Range("E16").Value = 7901

$ws.Range("B37").Value = 41193
$ws.Range("C37").Value = 2054
$ws.Range("D37").Value = 18859
$ws.Range("E37").Value = 20775
$ws.Range("G37").Value = 122
$ws.Range("H37").Value = 1559

$ws.Range("B38").Value = 41117
$ws.Range("C38").Value = 1109
$ws.Range("D38").Value = 18299
$ws.Range("E38").Value = 21732
$ws.Range("G38").Value = 19
$ws.Range("H38").Value = 1086

$ws.Range("B39").Value = 40866
$ws.Range("C39").Value = 451
$ws.Range("D39").Value = 26633
$ws.Range("E39").Value = 12678
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 1555

$ws.Range("B52").Value = 22800
$ws.Range("C52").Value = 400
$ws.Range("D52").Value = 16872
$ws.Range("E52").Value = 5614
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = 314

$ws.Range("B57").Value = 15834
$ws.Range("C57").Value = 361
$ws.Range("D57").Value = 11755
$ws.Range("E57").Value = 3976
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 103

$ws.Range("B58").Value = 15776
$ws.Range("C58").Value = 323
$ws.Range("D58").Value = 8765
$ws.Range("E58").Value = 6496
$ws.Range("G58").Value = 13
$ws.Range("H58").Value = 515

$ws.Range("B59").Value = 15619
$ws.Range("C59").Value = 800
$ws.Range("D59").Value = 2949
$ws.Range("E59").Value = 12047
$ws.Range("G59").Value = 22
$ws.Range("H59").Value = 623

$ws.Range("B63").Value = 12685
$ws.Range("C63").Value = 240
$ws.Range("D63").Value = 9066
$ws.Range("E63").Value = 2734
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 885

$ws.Range("B64").Value = 12675
$ws.Range("D64").Value = 11508
$ws.Range("E64").Value = 563
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 604

$ws.Range("B65").Value = 12602
$ws.Range("C65").Value = 39
$ws.Range("D65").Value = 11172
$ws.Range("E65").Value = 1148
$ws.Range("H65").Value = 282

$ws.Range("B66").Value = 12592
$ws.Range("D66").Value = 10100
$ws.Range("E66").Value = 2179
$ws.Range("H66").Value = 313

$ws.Range("B69").Value = 10923
$ws.Range("C69").Value = 53
$ws.Range("D69").Value = 7664
$ws.Range("E69").Value = 2913
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 346

$ws.Range("B70").Value = 9084
$ws.Range("C70").Value = 100
$ws.Range("D70").Value = 3912
$ws.Range("E70").Value = 4613
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 559

$ws.Range("B80").Value = 5747
$ws.Range("C80").Value = 56
$ws.Range("D80").Value = 4331
$ws.Range("E80").Value = 1364

$ws.Range("B91").Value = 4173
$ws.Range("C91").Value = 22
$ws.Range("E91").Value = 95

$ws.Range("B95").Value = 3343
$ws.Range("C95").Value = 22
$ws.Range("E95").Value = 1778

$ws.Range("B127").Value = 1104
$ws.Range("C127").Value = 18
$ws.Range("D127").Value = 830
$ws.Range("E127").Value = 265
$ws.Range("H127").Value = 9

$ws.Range("B128").Value = 1087
$ws.Range("D128").Value = 456
$ws.Range("E128").Value = 594
$ws.Range("H128").Value = 37

$ws.Range("B138").Value = 865
$ws.Range("C138").Value = 2
$ws.Range("D138").Value = 778
$ws.Range("E138").Value = 13

$ws.Range("B156").Value = 439
$ws.Range("C156").Value = 25
$ws.Range("E156").Value = 115

$ws.Range("D161").Value = 215
$ws.Range("E161").Value = 72

$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

$ws.Range("D213").Value = 10
$ws.Range("H213").Value = 1

